$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (ঢাকায়), shifting existing data down.
$ws.Rows("2:2").Insert()

# Excel's Insert() copies formatting down from the row above (the bold
# header row); strip that back to the plain/default style used by the
# other data rows before filling in values.
$ws.Range("A2:D2").ClearFormats()

# Populate the newly inserted row 2 with the "Unknown" location entry.
$ws.Range("A2").Value = "Unknown"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

# ঢাকায় (now row 3): normal count changes from 2 to 1.
$ws.Range("B3").Value = 1

# লক্ষ্মীপুরে (now row 4): severe count changes from 3 to 2.
$ws.Range("D4").Value = 2
